# Swap the content of rows 15 and 16 (katedra/jednotekPrednasek/jednotkaPrednasky)
# as described by the diff: row 15 and row 16 exchange their A, E and F values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "E", "F")

foreach ($col in $cols) {
    $cell15 = $ws.Range($col + "15")
    $cell16 = $ws.Range($col + "16")

    $v15 = $cell15.Value2
    $v16 = $cell16.Value2

    $cell15.Value2 = $v16
    $cell16.Value2 = $v15
}
